$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.039.11"
$ws.Range("E2").Value = "  -3.36%  "

$ws.Range("D3").Value = "3.326.66"
$ws.Range("E3").Value = "  -5.62%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.57%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "3.318.49"
$ws.Range("E9").Value = "  -5.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.620"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000273"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.59%  "

$ws.Range("D15").Value = "3.854.51"
$ws.Range("E15").Value = "  -5.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.11%  "

$ws.Range("E17").Value = "  -4.03%  "

$ws.Range("D18").Value = "3.320.95"
$ws.Range("E18").Value = "  -5.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.01%  "

$ws.Range("D20").Value = "63.955.36"
$ws.Range("E20").Value = "  -3.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.970"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "424.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.23%  "

$ws.Range("E24").Value = "  -2.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.53%  "

$ws.Range("E28").Value = "  -2.36%  "

$ws.Range("E29").Value = "  -5.85%  "

$ws.Range("E30").Value = "  -3.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "593.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.107"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.46%  "

$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.141"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.73%  "

$ws.Range("D40").Value = "0.0₃0748"
$ws.Range("E40").Value = "  -6.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.364"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.95%  "

$ws.Range("D42").Value = "3.093.63"
$ws.Range("E42").Value = "  -5.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("E44").Value = "  -5.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.10%  "

$ws.Range("E46").Value = "  -3.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.34%  "

$ws.Range("E48").Value = "  -3.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.35%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.94%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.95%  "
